# Week 9 time sheet updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task introduced: "Research/Personal Work" (row 8) ---
# Row 8 used to log its hour on Wednesday (D); it now logs on Saturday (G) and Sunday (H).
$ws.Range("A8").Value = "Research/Personal Work"
$ws.Range("D8").ClearContents()
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 1

# --- Row 9 ("Team Meting"): moved its hour from Monday (B) to Wednesday (D) ---
$ws.Range("A9").Value = "Team Meting"
$ws.Range("B9").ClearContents()
$ws.Range("D9").Value = 1

# --- Row 10 ("Sponsor Meeting"): logged a new hour on Monday (B) ---
$ws.Range("A10").Value = "Sponsor Meeting"
$ws.Range("B10").Value = 1

# --- Row 11 / Row 12 labels unaffected ---
$ws.Range("A11").Value = "TA Meeting"
$ws.Range("A12").Value = "Team Work"

# --- New row 13: "Organizing (misc.)" task, logged an hour on Saturday (G) ---
$ws.Range("A13").Value = "Organizing (misc.)"
$ws.Range("B13:F13").ClearContents()
$ws.Range("G13").Value = 1
$ws.Range("H13").ClearContents()
$ws.Range("I13").Formula = "=SUM(B13:H13)"

# --- Row 14: "Daily Total" row, now summing 6:13 instead of 6:12 ---
$ws.Range("A14").Value = "Daily Total"
$ws.Range("B14").Formula = "=SUM(B6:B13)"
$ws.Range("C14").Formula = "=SUM(C6:C13)"
$ws.Range("D14").Formula = "=SUM(D6:D13)"
$ws.Range("E14").Formula = "=SUM(E6:E13)"
$ws.Range("F14").Formula = "=SUM(F6:F13)"
$ws.Range("G14").Formula = "=SUM(G6:G13)"
$ws.Range("H14").Formula = "=SUM(H6:H13)"
$ws.Range("I14").Formula = "=SUM(I6:I13)"

# --- New row 22: repeats the closing note ---
$ws.Range("A22").Value = "It is intended both as an accountability tool and as validation for your estimates "

# --- Column A widened to fit the new, longer task labels ---
$ws.Columns("A").ColumnWidth = 19.33

# --- Selection moved to I15 ---
$ws.Range("I15").Select()
